$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.203683137893677
$ws.Range("B1").Value = 2.652856111526489
$ws.Range("C1").Value = 9.297165870666504
$ws.Range("D1").Value = 2.061556816101074
$ws.Range("E1").Value = 1.20168673992157
